$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 8 ---
# J8, N8 hold numeric/date-looking text; force Text format so Excel
# doesn't auto-convert them to a Number/Date serial.
$ws.Range("J8").NumberFormat = "@"
$ws.Range("J8").Value = "6233424802248"
$ws.Range("L8").Value = "Zrealizowana"
$ws.Range("N8").NumberFormat = "@"
$ws.Range("N8").Value = "2025-03-20"

# --- Add new row 9 ---
$ws.Range("A9").Value = 9
$ws.Range("B9").Value = "Biurko narożne Slash - Biały"
$ws.Range("C9").Value = "Kurier odmówił dostarczenia przesyłki, uszkodzona paczka."
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "19090"
$ws.Range("E9").Value = "Dorota "
$ws.Range("F9").Value = "Galinska"
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "+48726741823"
$ws.Range("H9").Value = "6wowjgbpc5+410544ae8@allegromail.pl"
$ws.Range("K9").Value = "Protokój ustalenia stanu przesyłki"
$ws.Range("L9").Value = "W trakcie"
$ws.Range("M9").NumberFormat = "@"
$ws.Range("M9").Value = "2025-03-20"

# --- Add new row 10 ---
$ws.Range("A10").Value = 10
$ws.Range("B10").Value = "Biurko narożne Slash - Beton Chicago jasnoszary"
$ws.Range("C10").Value = "Klient nie odebrał przesyłki. "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18936"
$ws.Range("E10").Value = "Mateusz "
$ws.Range("F10").Value = "Rakowski"
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "+48530394943"
$ws.Range("H10").Value = "1vbecjqnsb+27d32f7e1@allegromail.pl"
$ws.Range("K10").Value = "Wysłane ponownie, klient odebrał przesyłke."
$ws.Range("L10").Value = "Zrealizowana"
$ws.Range("M10").NumberFormat = "@"
$ws.Range("M10").Value = "2025-03-20"
$ws.Range("N10").NumberFormat = "@"
$ws.Range("N10").Value = "2025-03-20"
